# Keri.xlsx update: add FIXED/PARITY bit computation rows, move FC/CARD
# rows down by one, add references/links block, and restyle accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new blank row above row 2. This shifts the old rows 2-7
#    down to 3-8 and copies the formatting of the surrounding cells
#    (matches the style layout seen in the target file for row 2).
# ---------------------------------------------------------------------
$ws.Rows(2).Insert()

# New row 2 holds the bit-shift amounts used later by the FIXED/PARITY
# rows (K2=0, L2=1, M2=3, N2=31), plus a "FIXED" label in AQ2.
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 31
$ws.Range("AQ2").Value = "FIXED"

# ---------------------------------------------------------------------
# 2. Update the FC row (now row 5) with the new sample value.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "FC"
$ws.Range("B5").Value = 17

# ---------------------------------------------------------------------
# 3. Update the CARD row (now row 6) with the new sample value.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "CARD"
$ws.Range("B6").Value = 1337

# ---------------------------------------------------------------------
# 4. Row 7 computes the two fixed bits (bit 3 and bit 31).
# ---------------------------------------------------------------------
$ws.Range("M7").Formula = "=2^M2"
$ws.Range("N7").Formula = "=2^N2"
$ws.Range("AO7").Formula = "=SUM(H7:AM7)"
$ws.Range("AQ7").Value = "FIXED"

# ---------------------------------------------------------------------
# 5. Row 8 computes the parity bits from the hex digits of the running
#    total, then places them at bit 0 / bit 1.
# ---------------------------------------------------------------------
$ws.Range("D8").Formula = '=MOD(LEN(REGEXREPLACE(REGEXREPLACE(DEC2HEX(SUM(AO5:AO7)), "[236789CD]", "O"),"[^O]","")),2)'
$ws.Range("E8").Formula = '=MOD(LEN(REGEXREPLACE(REGEXREPLACE(DEC2HEX(SUM(AO5:AO7)), "[02578ADF]", "O"),"[^O]",""))+1,2)'
$ws.Range("K8").Formula = "=E8*2^K2"
$ws.Range("L8").Formula = "=D8*2^L2"
$ws.Range("AO8").Formula = "=SUM(H8:AM8)"
$ws.Range("AQ8").Value = "PARITY"

# ---------------------------------------------------------------------
# 6. Row 9 is the new TOTAL row (sum of the four contributing rows).
# ---------------------------------------------------------------------
$ws.Range("AO9").Formula = "=AO5+AO6+AO7+AO8"
$ws.Range("AP9").Formula = "=DEC2HEX(AO9, 8)"
$ws.Range("AQ9").Value = "TOTAL"

# ---------------------------------------------------------------------
# 7. Row 10 is the RESULT row, now pulling from AP9 instead of AP6.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "RESULT"
$ws.Range("B10").Formula = "=AP9"

# Clear out the old row-7 leftovers from the pre-insert layout (the
# RESULT row used to live at row 7 before the shift + edits above).
$ws.Range("A7").ClearContents()

# ---------------------------------------------------------------------
# 8. Reference / links block.
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "MAIN ALGO FROM"
$ws.Range("D13").Value = "https://github.com/Next-Flip/Momentum-Firmware/blob/dev/lib/lfrfid/protocols/protocol_keri.c"
$ws.Hyperlinks.Add($ws.Range("D13"), "https://github.com/Next-Flip/Momentum-Firmware/blob/dev/lib/lfrfid/protocols/protocol_keri.c") | Out-Null

$ws.Range("B14").Value = "FIXED/PARITY FROM"
$ws.Range("D14").Value = "https://github.com/RfidResearchGroup/proxmark3/blob/3ce68d4df918ef738686e7b63181dbe19809edd9/client/src/cmdlfkeri.c"
$ws.Hyperlinks.Add($ws.Range("D14"), "https://github.com/RfidResearchGroup/proxmark3/blob/3ce68d4df918ef738686e7b63181dbe19809edd9/client/src/cmdlfkeri.c") | Out-Null

$ws.Range("B16").Value = "ENTER YOUR FC AND CARD NUMBER IN CELLS B5 and B6. THE SPREADSHEET WILL CALCULATE THE HEX IN CELL B10."
$ws.Range("B16").Font.Color = 255

$ws.Range("B17").Value = "MORE INFO - "
$ws.Range("D17").Value = "https://github.com/jamisonderek/flipper-zero-tutorials/tree/main/rfid"
$ws.Hyperlinks.Add($ws.Range("D17"), "https://github.com/jamisonderek/flipper-zero-tutorials/tree/main/rfid") | Out-Null

# ---------------------------------------------------------------------
# 9. Conditional formatting ranges moved down with the FC/CARD rows.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
